# Refresh the crypto-price table to the latest scrape values.
# (Price/Volume columns and a handful of Coin/Link rows move as the
#  underlying coinranking.com ranking reshuffled between runs.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns hold numeric-looking values ("327.65",
# "-0.82%", ...) that are stored as plain text in the source sheet.
# Mark the cells we are about to rewrite as Text first so Excel/COM
# does not silently reinterpret them as numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '327.65'
$ws.Range("D3").Value = '44.36'
$ws.Range("E3").Value = '-0.82%'
$ws.Range("D4").Value = '5.292'
$ws.Range("E4").Value = '-4.65%'
$ws.Range("D5").Value = '0.08372'
$ws.Range("E5").Value = '1.96%'
$ws.Range("D6").Value = '1.929'
$ws.Range("E6").Value = '-6.11%'
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '0.9705'
$ws.Range("E7").Value = '-0.70%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '2.510'
$ws.Range("E8").Value = '-4.59%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1135'
$ws.Range("E9").Value = '1.32%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1901'
$ws.Range("E10").Value = '-0.42%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.09668'
$ws.Range("E11").Value = '-3.80%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.04600'
$ws.Range("E12").Value = '-2.28%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.1061'
$ws.Range("E13").Value = '0.21%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001298'
$ws.Range("E14").Value = '3.11%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.005792'
$ws.Range("E15").Value = '-2.16%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.387'
$ws.Range("E16").Value = '1.11%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.403'
$ws.Range("E17").Value = '-0.59%'
$ws.Range("E18").Value = '0.17%'
$ws.Range("D19").Value = '8.511'
$ws.Range("E19").Value = '-16.86%'
$ws.Range("D20").Value = '0.1388'
$ws.Range("E20").Value = '0.86%'
$ws.Range("E21").Value = '3.43%'
$ws.Range("D22").Value = '0.04159'
$ws.Range("E22").Value = '1.23%'
$ws.Range("D23").Value = '0.001232'
$ws.Range("E23").Value = '-5.30%'
$ws.Range("D24").Value = '0.004409'
$ws.Range("E24").Value = '0.44%'
$ws.Range("D25").Value = '0.0001300'
$ws.Range("E25").Value = '1.69%'
$ws.Range("D26").Value = '0.0002979'
$ws.Range("E26").Value = '-20.38%'
$ws.Range("E38").Value = '-1.87%'
$ws.Range("D39").Value = '0.05615'
$ws.Range("E39").Value = '-1.95%'
$ws.Range("D40").Value = '0.007821'
$ws.Range("E40").Value = '2.36%'
$ws.Range("D41").Value = '0.1414'
$ws.Range("E41").Value = '-0.79%'
$ws.Range("D42").Value = '0.007300'
$ws.Range("E42").Value = '-3.09%'
$ws.Range("D43").Value = '0.002050'
$ws.Range("E43").Value = '3.98%'
$ws.Range("D44").Value = '0.008671'
$ws.Range("E44").Value = '4.54%'
$ws.Range("D45").Value = '0.3514'
$ws.Range("D46").Value = '0.00006911'
$ws.Range("E46").Value = '-1.68%'
$ws.Range("E47").Value = '-0.04%'
$ws.Range("D48").Value = '0.003492'
$ws.Range("E48").Value = '-2.28%'
$ws.Range("D49").Value = '0.003530'
$ws.Range("E49").Value = '40.02%'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.04%'
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.04%'
